$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(15, 8).Value = 497.76
$ws.Cells.Item(15, 9).Value = 497.76
$ws.Cells.Item(15, 11).Value = 1493.28
$ws.Cells.Item(15, 13).Value = -1324.28
$ws.Cells.Item(19, 8).Value = 557.7778
$ws.Cells.Item(19, 9).Value = 653.5
$ws.Cells.Item(19, 11).Value = 653.5
$ws.Cells.Item(19, 13).Value = -478.5
$ws.Cells.Item(32, 8).Value = 2630
$ws.Cells.Item(32, 9).Value = 2250
$ws.Cells.Item(32, 10).Value = 2820
$ws.Cells.Item(32, 11).Value = 2250
$ws.Cells.Item(32, 12).Value = 2820
$ws.Cells.Item(32, 13).Value = -1924
$ws.Cells.Item(32, 14).Value = -3472
$ws.Cells.Item(33, 8).Value = 1016.8
$ws.Cells.Item(33, 9).Value = 174.82353
$ws.Cells.Item(33, 11).Value = 174.82353
$ws.Cells.Item(33, 13).Value = 54.17646999999999
$ws.Cells.Item(41, 8).Value = 1311.75
$ws.Cells.Item(41, 9).Value = 914.4
$ws.Cells.Item(41, 10).Value = 1974
$ws.Cells.Item(41, 11).Value = 914.4
$ws.Cells.Item(41, 12).Value = 1974
$ws.Cells.Item(41, 13).Value = -474.4
$ws.Cells.Item(41, 14).Value = -2854
$ws.Cells.Item(42, 8).Value = 58823972
$ws.Cells.Item(42, 9).Value = 100000140
$ws.Cells.Item(42, 10).Value = 869.4286
$ws.Cells.Item(42, 11).Value = 300000420
$ws.Cells.Item(42, 12).Value = 2608.2858
$ws.Cells.Item(42, 13).Value = -300000190
$ws.Cells.Item(42, 14).Value = -3068.2858
$ws.Cells.Item(55, 8).Value = 147.6875
$ws.Cells.Item(55, 9).Value = 88.333336
$ws.Cells.Item(55, 10).Value = 325.75
$ws.Cells.Item(55, 11).Value = 88.333336
$ws.Cells.Item(55, 12).Value = 325.75
$ws.Cells.Item(55, 13).Value = 125.666664
$ws.Cells.Item(55, 14).Value = -753.75
$ws.Cells.Item(74, 8).Value = 13721.75
$ws.Cells.Item(74, 9).Value = 9962.333000000001
$ws.Cells.Item(74, 11).Value = 9962.333000000001
$ws.Cells.Item(74, 13).Value = -9026.333000000001
$ws.Cells.Item(76, 8).Value = 3000
$ws.Cells.Item(76, 10).Value = 3000
$ws.Cells.Item(76, 12).Value = 3000
$ws.Cells.Item(76, 14).Value = -3630
$ws.Cells.Item(77, 8).Value = 13721.75
$ws.Cells.Item(77, 9).Value = 9962.333000000001
$ws.Cells.Item(77, 11).Value = 49811.665
$ws.Cells.Item(77, 13).Value = -45131.665
$ws.Cells.Item(79, 8).Value = 3000
$ws.Cells.Item(79, 10).Value = 3000
$ws.Cells.Item(79, 12).Value = 3000
$ws.Cells.Item(79, 14).Value = -5184
$ws.Cells.Item(92, 8).Value = 632.8095
$ws.Cells.Item(92, 9).Value = 701.3889
$ws.Cells.Item(92, 11).Value = 701.3889
$ws.Cells.Item(92, 13).Value = 546.6111
$ws.Cells.Item(98, 8).Value = 2810.75
$ws.Cells.Item(98, 9).Value = 2508.6428
$ws.Cells.Item(98, 10).Value = 3515.6667
$ws.Cells.Item(98, 11).Value = 2508.6428
$ws.Cells.Item(98, 12).Value = 3515.6667
$ws.Cells.Item(98, 13).Value = -1010.6428
$ws.Cells.Item(98, 14).Value = -6511.6667
$ws.Cells.Item(99, 8).Value = 403.16666
$ws.Cells.Item(99, 9).Value = 403.16666
$ws.Cells.Item(99, 11).Value = 1209.49998
$ws.Cells.Item(99, 13).Value = 288.5000199999999
$ws.Cells.Item(106, 8).Value = 4344
$ws.Cells.Item(106, 9).Value = 4142.6665
$ws.Cells.Item(106, 11).Value = 4142.6665
$ws.Cells.Item(106, 13).Value = -3511.6665
$ws.Cells.Item(107, 8).Value = 95352.5
$ws.Cells.Item(107, 9).Value = 95352.5
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 95352.5
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -93432.5
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(112, 8).Value = 3995.5454
$ws.Cells.Item(112, 10).Value = 3995.5454
$ws.Cells.Item(112, 12).Value = 11986.6362
$ws.Cells.Item(112, 14).Value = -14202.6362
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(122, 8).Value = 2810.75
$ws.Cells.Item(122, 9).Value = 2508.6428
$ws.Cells.Item(122, 10).Value = 3515.6667
$ws.Cells.Item(122, 11).Value = 7525.928400000001
$ws.Cells.Item(122, 12).Value = 10547.0001
$ws.Cells.Item(122, 13).Value = -5075.928400000001
$ws.Cells.Item(122, 14).Value = -15447.0001
$ws.Cells.Item(131, 8).Value = 7656.143
$ws.Cells.Item(131, 9).Value = 5864.6665
$ws.Cells.Item(131, 10).Value = 8999.75
$ws.Cells.Item(131, 11).Value = 17593.9995
$ws.Cells.Item(131, 12).Value = 26999.25
$ws.Cells.Item(131, 13).Value = -12553.9995
$ws.Cells.Item(131, 14).Value = -37079.25
$ws.Cells.Item(132, 8).Value = 2823.1333
$ws.Cells.Item(132, 9).Value = 2310.0715
$ws.Cells.Item(132, 11).Value = 6930.2145
$ws.Cells.Item(132, 13).Value = -4400.2145
$ws.Cells.Item(136, 8).Value = 190337.25
$ws.Cells.Item(136, 10).Value = 190337.25
$ws.Cells.Item(136, 12).Value = 190337.25
$ws.Cells.Item(136, 14).Value = -200537.25
$ws.Cells.Item(137, 8).Value = 6787.3335
$ws.Cells.Item(137, 9).Value = 8157.684
$ws.Cells.Item(137, 10).Value = 1580
$ws.Cells.Item(137, 11).Value = 24473.052
$ws.Cells.Item(137, 12).Value = 4740
$ws.Cells.Item(137, 13).Value = -21923.052
$ws.Cells.Item(137, 14).Value = -9840
$ws.Cells.Item(138, 8).Value = 4565
$ws.Cells.Item(138, 9).Value = 3666.5
$ws.Cells.Item(138, 10).Value = 4901.9375
$ws.Cells.Item(138, 11).Value = 10999.5
$ws.Cells.Item(138, 12).Value = 14705.8125
$ws.Cells.Item(138, 13).Value = -5859.5
$ws.Cells.Item(138, 14).Value = -24985.8125
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
$ws.Cells.Item(141, 8).Value = 4773.409
$ws.Cells.Item(141, 9).Value = 4811.5
$ws.Cells.Item(141, 10).Value = 4741.6665
$ws.Cells.Item(141, 11).Value = 14434.5
$ws.Cells.Item(141, 12).Value = 14224.9995
$ws.Cells.Item(141, 13).Value = -9254.5
$ws.Cells.Item(141, 14).Value = -24584.9995

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 2605
$ws.Cells.Item(2, 9).Value = 912.4
$ws.Cells.Item(2, 10).Value = 3545.3333
$ws.Cells.Item(2, 11).Value = 912.4
$ws.Cells.Item(2, 12).Value = 3545.3333
$ws.Cells.Item(2, 13).Value = -799.4
$ws.Cells.Item(2, 14).Value = -3771.3333
$ws.Cells.Item(4, 8).Value = 330
$ws.Cells.Item(4, 9).Value = 50
$ws.Cells.Item(4, 10).Value = 400
$ws.Cells.Item(4, 11).Value = 50
$ws.Cells.Item(4, 12).Value = 400
$ws.Cells.Item(4, 13).Value = 66
$ws.Cells.Item(4, 14).Value = -632
$ws.Cells.Item(13, 8).Value = 3460950.2
$ws.Cells.Item(13, 10).Value = 191175
$ws.Cells.Item(13, 12).Value = 191175
$ws.Cells.Item(13, 14).Value = -191463
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 14).ClearContents()
$ws.Cells.Item(32, 8).Value = 3711
$ws.Cells.Item(32, 9).Value = 3736.9375
$ws.Cells.Item(32, 11).Value = 3736.9375
$ws.Cells.Item(32, 13).Value = -3449.9375
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).ClearContents()
$ws.Cells.Item(45, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 3945.9092
$ws.Cells.Item(61, 9).Value = 3871.2415
$ws.Cells.Item(61, 11).Value = 3871.2415
$ws.Cells.Item(61, 13).Value = -3659.2415
$ws.Cells.Item(63, 8).Value = 7166.6665
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 7166.6665
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 7166.6665
$ws.Cells.Item(63, 13).ClearContents()
$ws.Cells.Item(63, 14).Value = -8538.666499999999
$ws.Cells.Item(66, 8).Value = 7166.6665
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 7166.6665
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 35833.3325
$ws.Cells.Item(66, 13).ClearContents()
$ws.Cells.Item(66, 14).Value = -42697.3325
$ws.Cells.Item(74, 8).Value = 3965.5
$ws.Cells.Item(74, 9).Value = 4191.273
$ws.Cells.Item(74, 10).Value = 3137.6667
$ws.Cells.Item(74, 11).Value = 4191.273
$ws.Cells.Item(74, 12).Value = 3137.6667
$ws.Cells.Item(74, 13).Value = -3317.273
$ws.Cells.Item(74, 14).Value = -4885.6667
$ws.Cells.Item(77, 8).Value = 3965.5
$ws.Cells.Item(77, 9).Value = 4191.273
$ws.Cells.Item(77, 10).Value = 3137.6667
$ws.Cells.Item(77, 11).Value = 20956.365
$ws.Cells.Item(77, 12).Value = 15688.3335
$ws.Cells.Item(77, 13).Value = -16588.365
$ws.Cells.Item(77, 14).Value = -24424.3335
$ws.Cells.Item(97, 8).Value = 323.4375
$ws.Cells.Item(97, 9).Value = 317.6
$ws.Cells.Item(97, 10).Value = 411
$ws.Cells.Item(97, 11).Value = 317.6
$ws.Cells.Item(97, 12).Value = 411
$ws.Cells.Item(97, 13).Value = 178.4
$ws.Cells.Item(97, 14).Value = -1403
$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 14).ClearContents()
$ws.Cells.Item(102, 8).Value = 3649.1667
$ws.Cells.Item(102, 9).Value = 1834.6428
$ws.Cells.Item(102, 10).Value = 10000
$ws.Cells.Item(102, 11).Value = 1834.6428
$ws.Cells.Item(102, 12).Value = 10000
$ws.Cells.Item(102, 13).Value = -212.6428000000001
$ws.Cells.Item(102, 14).Value = -13244
$ws.Cells.Item(110, 8).Value = 836.7778
$ws.Cells.Item(110, 9).Value = 836.7778
$ws.Cells.Item(110, 11).Value = 836.7778
$ws.Cells.Item(110, 13).Value = 1208.2222
$ws.Cells.Item(116, 8).Value = 2605
$ws.Cells.Item(116, 9).Value = 912.4
$ws.Cells.Item(116, 10).Value = 3545.3333
$ws.Cells.Item(116, 11).Value = 912.4
$ws.Cells.Item(116, 12).Value = 3545.3333
$ws.Cells.Item(116, 13).Value = 1381.6
$ws.Cells.Item(116, 14).Value = -8133.3333
$ws.Cells.Item(120, 8).Value = 85000
$ws.Cells.Item(120, 10).Value = 85000
$ws.Cells.Item(120, 12).Value = 85000
$ws.Cells.Item(120, 14).Value = -94676
$ws.Cells.Item(121, 8).Value = 114999.5
$ws.Cells.Item(121, 10).Value = 114999.5
$ws.Cells.Item(121, 12).Value = 114999.5
$ws.Cells.Item(121, 14).Value = -118493.5
$ws.Cells.Item(122, 8).Value = 13892522
$ws.Cells.Item(122, 9).Value = 15876526
$ws.Cells.Item(122, 10).Value = 4500
$ws.Cells.Item(122, 11).Value = 47629578
$ws.Cells.Item(122, 12).Value = 13500
$ws.Cells.Item(122, 13).Value = -47627128
$ws.Cells.Item(122, 14).Value = -18400
$ws.Cells.Item(132, 8).Value = 3281.9614
$ws.Cells.Item(132, 9).Value = 3315.9167
$ws.Cells.Item(132, 11).Value = 9947.750100000001
$ws.Cells.Item(132, 13).Value = -7417.750100000001
$ws.Cells.Item(136, 8).Value = 3945.9092
$ws.Cells.Item(136, 9).Value = 3871.2415
$ws.Cells.Item(136, 11).Value = 11613.7245
$ws.Cells.Item(136, 13).Value = -9063.7245

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 2605
$ws.Cells.Item(3, 9).Value = 912.4
$ws.Cells.Item(3, 10).Value = 3545.3333
$ws.Cells.Item(3, 11).Value = 912.4
$ws.Cells.Item(3, 12).Value = 3545.3333
$ws.Cells.Item(3, 13).Value = -798.4
$ws.Cells.Item(3, 14).Value = -3773.3333
$ws.Cells.Item(11, 8).Value = 1078.9231
$ws.Cells.Item(11, 10).Value = 1876.25
$ws.Cells.Item(11, 12).Value = 1876.25
$ws.Cells.Item(11, 14).Value = -2156.25
$ws.Cells.Item(22, 8).Value = 1142.0714
$ws.Cells.Item(22, 9).Value = 1018.75
$ws.Cells.Item(22, 11).Value = 1018.75
$ws.Cells.Item(22, 13).Value = -845.75
$ws.Cells.Item(86, 8).Value = 6702.778
$ws.Cells.Item(86, 9).Value = 4332.143
$ws.Cells.Item(86, 11).Value = 4332.143
$ws.Cells.Item(86, 13).Value = -3209.143
$ws.Cells.Item(89, 8).Value = 6702.778
$ws.Cells.Item(89, 9).Value = 4332.143
$ws.Cells.Item(89, 11).Value = 21660.715
$ws.Cells.Item(89, 13).Value = -16044.715
$ws.Cells.Item(99, 8).Value = 4410.6
$ws.Cells.Item(99, 9).Value = 2800.8572
$ws.Cells.Item(99, 11).Value = 2800.8572
$ws.Cells.Item(99, 13).Value = -1302.8572
$ws.Cells.Item(105, 8).Value = 1614.625
$ws.Cells.Item(105, 9).Value = 1137.9333
$ws.Cells.Item(105, 10).Value = 2409.111
$ws.Cells.Item(105, 11).Value = 1137.9333
$ws.Cells.Item(105, 12).Value = 2409.111
$ws.Cells.Item(105, 13).Value = 609.0667000000001
$ws.Cells.Item(105, 14).Value = -5903.111
$ws.Cells.Item(107, 8).Value = 4617.1763
$ws.Cells.Item(107, 10).Value = 8124.125
$ws.Cells.Item(107, 12).Value = 8124.125
$ws.Cells.Item(107, 14).Value = -11964.125
$ws.Cells.Item(137, 8).Value = 75000
$ws.Cells.Item(137, 10).Value = 75000
$ws.Cells.Item(137, 12).Value = 75000
$ws.Cells.Item(137, 14).Value = -85200

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 8).Value = 3000
$ws.Cells.Item(4, 10).Value = 3000
$ws.Cells.Item(4, 12).Value = 3000
$ws.Cells.Item(4, 14).Value = -3224
$ws.Cells.Item(31, 8).Value = 4208.095
$ws.Cells.Item(31, 9).Value = 1862.3334
$ws.Cells.Item(31, 11).Value = 1862.3334
$ws.Cells.Item(31, 13).Value = -1567.3334
$ws.Cells.Item(34, 8).Value = 4208.095
$ws.Cells.Item(34, 9).Value = 1862.3334
$ws.Cells.Item(34, 11).Value = 1862.3334
$ws.Cells.Item(34, 13).Value = -1660.3334
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).ClearContents()
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 3650
$ws.Cells.Item(86, 9).Value = 3800
$ws.Cells.Item(86, 11).Value = 3800
$ws.Cells.Item(86, 13).Value = -2677
$ws.Cells.Item(89, 8).Value = 3650
$ws.Cells.Item(89, 9).Value = 3800
$ws.Cells.Item(89, 11).Value = 19000
$ws.Cells.Item(89, 13).Value = -13384
$ws.Cells.Item(93, 8).Value = 32598
$ws.Cells.Item(93, 9).Value = 23630.666
$ws.Cells.Item(93, 10).Value = 59500
$ws.Cells.Item(93, 11).Value = 23630.666
$ws.Cells.Item(93, 12).Value = 59500
$ws.Cells.Item(93, 13).Value = -21758.666
$ws.Cells.Item(93, 14).Value = -63244
$ws.Cells.Item(99, 8).Value = 2354.9375
$ws.Cells.Item(99, 9).Value = 2404.6428
$ws.Cells.Item(99, 11).Value = 2404.6428
$ws.Cells.Item(99, 13).Value = -906.6428000000001
$ws.Cells.Item(104, 8).Value = 30159.5
$ws.Cells.Item(104, 9).Value = 30159.5
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 11).Value = 30159.5
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 13).Value = -27538.5
$ws.Cells.Item(104, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 1841.6316
$ws.Cells.Item(107, 9).Value = 997
$ws.Cells.Item(107, 10).Value = 2000
$ws.Cells.Item(107, 11).Value = 997
$ws.Cells.Item(107, 12).Value = 2000
$ws.Cells.Item(107, 13).Value = 923
$ws.Cells.Item(107, 14).Value = -5840
$ws.Cells.Item(122, 8).Value = 1874.75
$ws.Cells.Item(122, 9).Value = 1874.75
$ws.Cells.Item(122, 11).Value = 5624.25
$ws.Cells.Item(122, 13).Value = -3174.25
$ws.Cells.Item(126, 8).Value = 2354.9375
$ws.Cells.Item(126, 9).Value = 2404.6428
$ws.Cells.Item(126, 11).Value = 7213.928400000001
$ws.Cells.Item(126, 13).Value = -4743.928400000001
$ws.Cells.Item(132, 8).Value = 1999.3125
$ws.Cells.Item(132, 9).Value = 1999.3125
$ws.Cells.Item(132, 11).Value = 5997.9375
$ws.Cells.Item(132, 13).Value = -3467.9375
$ws.Cells.Item(134, 8).Value = 1962.1578
$ws.Cells.Item(134, 9).Value = 2043.7778
$ws.Cells.Item(134, 10).Value = 493
$ws.Cells.Item(134, 11).Value = 6131.3334
$ws.Cells.Item(134, 12).Value = 1479
$ws.Cells.Item(134, 13).Value = -3596.3334
$ws.Cells.Item(134, 14).Value = -6549

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(23, 8).Value = 5454.3335
$ws.Cells.Item(23, 9).Value = 8254.333000000001
$ws.Cells.Item(23, 10).Value = 3587.6667
$ws.Cells.Item(23, 11).Value = 24762.999
$ws.Cells.Item(23, 12).Value = 10763.0001
$ws.Cells.Item(23, 13).Value = -24527.999
$ws.Cells.Item(23, 14).Value = -11233.0001
$ws.Cells.Item(29, 8).Value = 725
$ws.Cells.Item(29, 9).Value = 633.6667
$ws.Cells.Item(29, 10).Value = 999
$ws.Cells.Item(29, 11).Value = 1901.0001
$ws.Cells.Item(29, 12).Value = 2997
$ws.Cells.Item(29, 13).Value = -1624.0001
$ws.Cells.Item(29, 14).Value = -3551
$ws.Cells.Item(47, 8).Value = 497.66666
$ws.Cells.Item(47, 9).Value = 497.66666
$ws.Cells.Item(47, 11).Value = 1492.99998
$ws.Cells.Item(47, 13).Value = -1061.99998
$ws.Cells.Item(50, 8).Value = 1028
$ws.Cells.Item(50, 9).Value = 925.55554
$ws.Cells.Item(50, 11).Value = 2776.66662
$ws.Cells.Item(50, 13).Value = -2295.66662
$ws.Cells.Item(53, 8).Value = 1028
$ws.Cells.Item(53, 9).Value = 925.55554
$ws.Cells.Item(53, 11).Value = 2776.66662
$ws.Cells.Item(53, 13).Value = -2295.66662
$ws.Cells.Item(99, 8).Value = 12792.182
$ws.Cells.Item(99, 9).Value = 2579.3333
$ws.Cells.Item(99, 10).Value = 58750
$ws.Cells.Item(99, 11).Value = 7737.999899999999
$ws.Cells.Item(99, 12).Value = 176250
$ws.Cells.Item(99, 13).Value = -5491.999899999999
$ws.Cells.Item(99, 14).Value = -180742
$ws.Cells.Item(104, 8).Value = 26322.666
$ws.Cells.Item(104, 10).Value = 60000
$ws.Cells.Item(104, 12).Value = 180000
$ws.Cells.Item(104, 14).Value = -185242
$ws.Cells.Item(107, 8).Value = 3426.2703
$ws.Cells.Item(107, 9).Value = 5260.2
$ws.Cells.Item(107, 10).Value = 2747.037
$ws.Cells.Item(107, 11).Value = 15780.6
$ws.Cells.Item(107, 12).Value = 8241.110999999999
$ws.Cells.Item(107, 13).Value = -13860.6
$ws.Cells.Item(107, 14).Value = -12081.111
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).ClearContents()
$ws.Cells.Item(136, 14).ClearContents()

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(15, 8).Value = 23969
$ws.Cells.Item(15, 10).Value = 23969
$ws.Cells.Item(15, 12).Value = 23969
$ws.Cells.Item(15, 14).Value = -24545
$ws.Cells.Item(22, 8).Value = 12333.333
$ws.Cells.Item(22, 9).Value = 8500
$ws.Cells.Item(22, 10).Value = 20000
$ws.Cells.Item(22, 11).Value = 8500
$ws.Cells.Item(22, 12).Value = 20000
$ws.Cells.Item(22, 13).Value = -7971
$ws.Cells.Item(22, 14).Value = -21058
$ws.Cells.Item(70, 8).Value = 16675557
$ws.Cells.Item(70, 9).Value = 111114080
$ws.Cells.Item(70, 11).Value = 111114080
$ws.Cells.Item(70, 13).Value = -111113810
$ws.Cells.Item(73, 8).Value = 16675557
$ws.Cells.Item(73, 9).Value = 111114080
$ws.Cells.Item(73, 11).Value = 111114080
$ws.Cells.Item(73, 13).Value = -111113144
$ws.Cells.Item(80, 8).Value = 10036.117
$ws.Cells.Item(80, 9).Value = 13646
$ws.Cells.Item(80, 11).Value = 13646
$ws.Cells.Item(80, 13).Value = -12648
$ws.Cells.Item(81, 8).Value = 23969
$ws.Cells.Item(81, 10).Value = 23969
$ws.Cells.Item(81, 12).Value = 23969
$ws.Cells.Item(81, 14).Value = -25965
$ws.Cells.Item(83, 8).Value = 10036.117
$ws.Cells.Item(83, 9).Value = 13646
$ws.Cells.Item(83, 11).Value = 68230
$ws.Cells.Item(83, 13).Value = -63238
$ws.Cells.Item(84, 8).Value = 23969
$ws.Cells.Item(84, 10).Value = 23969
$ws.Cells.Item(84, 12).Value = 71907
$ws.Cells.Item(84, 14).Value = -81891
$ws.Cells.Item(97, 8).Value = 4312.7744
$ws.Cells.Item(97, 9).Value = 1239.6923
$ws.Cells.Item(97, 11).Value = 1239.6923
$ws.Cells.Item(97, 13).Value = -743.6922999999999
$ws.Cells.Item(102, 8).Value = 2438.2307
$ws.Cells.Item(102, 9).Value = 2438.2307
$ws.Cells.Item(102, 11).Value = 2438.2307
$ws.Cells.Item(102, 13).Value = -816.2307000000001
$ws.Cells.Item(126, 8).Value = 5922.222
$ws.Cells.Item(126, 9).Value = 6100.8
$ws.Cells.Item(126, 11).Value = 18302.4
$ws.Cells.Item(126, 13).Value = -15832.4
$ws.Cells.Item(132, 8).Value = 3067.5
$ws.Cells.Item(132, 9).Value = 2869.5454
$ws.Cells.Item(132, 10).Value = 5245
$ws.Cells.Item(132, 11).Value = 8608.636200000001
$ws.Cells.Item(132, 12).Value = 15735
$ws.Cells.Item(132, 13).Value = -6078.636200000001
$ws.Cells.Item(132, 14).Value = -20795

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 8).Value = 1200
$ws.Cells.Item(2, 9).Value = 1200
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 1200
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -1088
$ws.Cells.Item(2, 14).ClearContents()
$ws.Cells.Item(12, 8).Value = 3611.7144
$ws.Cells.Item(12, 9).Value = 499
$ws.Cells.Item(12, 10).Value = 4130.5
$ws.Cells.Item(12, 11).Value = 499
$ws.Cells.Item(12, 12).Value = 4130.5
$ws.Cells.Item(12, 13).Value = -329
$ws.Cells.Item(12, 14).Value = -4470.5
$ws.Cells.Item(40, 8).Value = 5876.027
$ws.Cells.Item(40, 9).Value = 3494.25
$ws.Cells.Item(40, 10).Value = 7690.7144
$ws.Cells.Item(40, 11).Value = 3494.25
$ws.Cells.Item(40, 12).Value = 7690.7144
$ws.Cells.Item(40, 13).Value = -3358.25
$ws.Cells.Item(40, 14).Value = -7962.7144
$ws.Cells.Item(55, 8).Value = 1098
$ws.Cells.Item(55, 9).Value = 172.33333
$ws.Cells.Item(55, 11).Value = 172.33333
$ws.Cells.Item(55, 13).Value = 0.6666700000000105
$ws.Cells.Item(61, 8).Value = 4227.7666
$ws.Cells.Item(61, 9).Value = 2384.125
$ws.Cells.Item(61, 11).Value = 2384.125
$ws.Cells.Item(61, 13).Value = -2182.125
$ws.Cells.Item(82, 8).Value = 2827.3572
$ws.Cells.Item(82, 9).Value = 2526.3635
$ws.Cells.Item(82, 10).Value = 3931
$ws.Cells.Item(82, 11).Value = 2526.3635
$ws.Cells.Item(82, 12).Value = 3931
$ws.Cells.Item(82, 13).Value = -2165.3635
$ws.Cells.Item(82, 14).Value = -4653
$ws.Cells.Item(85, 8).Value = 2827.3572
$ws.Cells.Item(85, 9).Value = 2526.3635
$ws.Cells.Item(85, 10).Value = 3931
$ws.Cells.Item(85, 11).Value = 2526.3635
$ws.Cells.Item(85, 12).Value = 3931
$ws.Cells.Item(85, 13).Value = -1278.3635
$ws.Cells.Item(85, 14).Value = -6427
$ws.Cells.Item(107, 8).Value = 6659.6665
$ws.Cells.Item(107, 9).Value = 6659.6665
$ws.Cells.Item(107, 11).Value = 6659.6665
$ws.Cells.Item(107, 13).Value = -4739.6665
$ws.Cells.Item(113, 8).Value = 4227.7666
$ws.Cells.Item(113, 9).Value = 2384.125
$ws.Cells.Item(113, 11).Value = 2384.125
$ws.Cells.Item(113, 13).Value = -214.125
$ws.Cells.Item(132, 8).Value = 5407.654
$ws.Cells.Item(132, 9).Value = 6137.375
$ws.Cells.Item(132, 11).Value = 18412.125
$ws.Cells.Item(132, 13).Value = -15882.125

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(5, 8).Value = 50019950
$ws.Cells.Item(5, 9).Value = 39800
$ws.Cells.Item(5, 11).Value = 39800
$ws.Cells.Item(5, 13).Value = -39688
$ws.Cells.Item(62, 8).Value = 11950
$ws.Cells.Item(65, 8).Value = 11950
$ws.Cells.Item(93, 8).Value = 44999.5
$ws.Cells.Item(93, 10).Value = 44999.5
$ws.Cells.Item(93, 12).Value = 44999.5
$ws.Cells.Item(93, 14).Value = -49991.5
$ws.Cells.Item(107, 8).Value = 6599.6
$ws.Cells.Item(107, 9).Value = 6333
$ws.Cells.Item(107, 10).Value = 6999.5
$ws.Cells.Item(107, 11).Value = 18999
$ws.Cells.Item(107, 12).Value = 20998.5
$ws.Cells.Item(107, 13).Value = -17079
$ws.Cells.Item(107, 14).Value = -24838.5
$ws.Cells.Item(113, 8).Value = 915.8421
$ws.Cells.Item(113, 9).Value = 794.3
$ws.Cells.Item(113, 10).Value = 1371.625
$ws.Cells.Item(113, 11).Value = 2382.9
$ws.Cells.Item(113, 12).Value = 4114.875
$ws.Cells.Item(113, 13).Value = -212.8999999999996
$ws.Cells.Item(113, 14).Value = -8454.875
$ws.Cells.Item(122, 8).Value = 5092.385
$ws.Cells.Item(122, 9).Value = 3006.1177
$ws.Cells.Item(122, 11).Value = 9018.3531
$ws.Cells.Item(122, 13).Value = -6568.3531
$ws.Cells.Item(126, 8).Value = 3162
$ws.Cells.Item(126, 9).Value = 3173.5715
$ws.Cells.Item(126, 11).Value = 9520.7145
$ws.Cells.Item(126, 13).Value = -7050.7145
$ws.Cells.Item(132, 8).Value = 10473.5625
$ws.Cells.Item(132, 9).Value = 6346.75
$ws.Cells.Item(132, 10).Value = 11849.167
$ws.Cells.Item(132, 11).Value = 19040.25
$ws.Cells.Item(132, 12).Value = 35547.501
$ws.Cells.Item(132, 13).Value = -16510.25
$ws.Cells.Item(132, 14).Value = -40607.501
$ws.Cells.Item(136, 8).Value = 2769.2222
$ws.Cells.Item(136, 9).Value = 2089.4
$ws.Cells.Item(136, 10).Value = 6168.3335
$ws.Cells.Item(136, 11).Value = 6268.200000000001
$ws.Cells.Item(136, 12).Value = 18505.0005
$ws.Cells.Item(136, 13).Value = -3718.200000000001
$ws.Cells.Item(136, 14).Value = -23605.0005
